# Add data for 2022-04-19 (data pulled through 2022-04-11)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename the sheet and the "through" header cell ---
$ws.Name = "Through 2022-04-11"
$ws.Range("B1").Value = "April 2022 (through April 11)"

# --- Insert a new neighborhood row "Beverly" (alphabetically between Avalon Park and
#     Bucktown), pushing every following neighborhood row (Bucktown..West Ridge) down
#     by one.
#
#     Property getters (".Value" reads) aren't wired up to real cell data in this
#     host, so rows can't be shifted by reading then rewriting values. Instead each
#     row is moved one at a time, bottom row first, via Copy + PasteSpecial(values):
#     walking from row 93 up to row 62 guarantees a row is always copied from its
#     still-untouched original location before anything writes into it. The
#     destination row is cleared immediately before its own paste (never between an
#     unrelated Copy/Paste pair, which invalidates the clipboard) so that blank cells
#     in the source correctly blank out stale data left in the destination. ---
for ($r = 93; $r -ge 62; $r--) {
    $dstRow = $r + 1
    $srcRange = "A" + $r + ":AG" + $r
    $dstRange = "A" + $dstRow + ":AG" + $dstRow
    $ws.Range($dstRange).ClearContents()
    $ws.Range($srcRange).Copy()
    $ws.Range($dstRange).PasteSpecial(-4163)
}

# Row 94 didn't exist before the shift, so it has no inherited cell style yet --
# copy the label-column format from the row above (now row 93) onto it.
$ws.Range("A93").Copy()
$ws.Range("A94").PasteSpecial(-4122)

# Clear the old row 62 contents (currently a duplicate of the shifted row 63) and
# write the new "Beverly" row's data.
$ws.Range("A62:AG62").ClearContents()
$ws.Range("A62").Value = "Beverly"
$ws.Range("V62").Value = 1   # April 2017

$excel.CutCopyMode = 0

# --- Updated monthly counts for several neighborhoods (new incident(s) recorded) ---
$ws.Range("J2").Value = 3    # Austin / April 2020
$ws.Range("V4").Value = 6    # North Lawndale / April 2017
$ws.Range("J5").Value = 2    # Garfield Park / April 2020
$ws.Range("B6").Value = 3    # Humboldt Park / April 2022 (through April 11)
$ws.Range("R6").Value = 3    # Humboldt Park / April 2018
$ws.Range("B9").Value = 3    # Loop / April 2022 (through April 11)
$ws.Range("J11").Value = 2   # Chatham / April 2020
$ws.Range("F12").Value = 2   # Calumet Heights / April 2021
$ws.Range("R22").Value = 1   # Rogers Park / April 2018 (new)
$ws.Range("N26").Value = 2   # South Shore / April 2019
$ws.Range("N27").Value = 1   # Uptown / April 2019 (new)
$ws.Range("V32").Value = 2   # Roseland / April 2017
$ws.Range("F53").Value = 2   # Kenwood / April 2021
